# Replace the hard-coded "2024" year in the certificate date line with a
# dynamic {year} placeholder, split into its own run (matching the diff,
# which shows the original single run becoming two runs: "of {month}, "
# and "{year}").

$d = $word.ActiveDocument

# Locate "of {month}, 2024" and replace it with "of {month}, " (keeps the
# original run's formatting: sz=22 / szCs=22). Find/Replace preserves the
# run's rPr and Word automatically marks the trailing space as
# xml:space="preserve".
$r = $d.Content
$found = $r.Find.Execute("of {month}, 2024", $true, $false, $false, $false, $false, $true, 1, $false, "of {month}, ", 2)

if (-not $found) {
    throw "Could not find 'of {month}, 2024' in the document"
}

# $r now spans the just-replaced "of {month}, " text. Remember where it
# ends -- that's where the new {year} run needs to be inserted.
$insertPos = $r.End
$r.Collapse(0)
$r.InsertAfter("{year}")

# Word's "type at cursor" semantics merge the freshly inserted text into
# the preceding run when the formatting matches exactly. Touch the new
# range's Bold property (set then immediately clear it) so it is forced
# to stay in its own run, matching the target run-split in the diff, while
# leaving its effective/visible formatting unchanged.
$yearRange = $d.Range($insertPos, $insertPos + 6)
$yearRange.Bold = 1
$yearRange.Bold = 0
